$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Coliflor, Macroferia Regional de Talca) was added
# as the new row 203; every subsequent record (old rows 203-256) shifts down
# by one row to make room, so the sheet grows from A1:R256 to A1:R257.
$ws.Rows.Item(203).Insert()

$ws.Cells.Item(203, 1).Value = 5
$ws.Cells.Item(203, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(203, 3).Value = "Maule"
$ws.Cells.Item(203, 4).Value = 44736
$ws.Cells.Item(203, 5).Value = 7
$ws.Cells.Item(203, 6).Value = 100112008
$ws.Cells.Item(203, 7).Value = "Coliflor"
$ws.Cells.Item(203, 8).Value = "Sin especificar"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 3000
$ws.Cells.Item(203, 11).Value = 1000
$ws.Cells.Item(203, 12).Value = 1000
$ws.Cells.Item(203, 13).Value = 1000
$ws.Cells.Item(203, 14).Value = "$/unidad"
$ws.Cells.Item(203, 15).Value = "Región del Maule"
$ws.Cells.Item(203, 16).Value = 1000
$ws.Cells.Item(203, 17).Value = 1
$ws.Cells.Item(203, 18).Value = "Hortaliza"
